# Generate Report for Handback
# Update timestamps and priority values produced by a new handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for rows 4-5
$wsOverview.Range("G4").Value = "2016-08-22 18:16:28"
$wsOverview.Range("G5").Value = "2016-08-22 18:16:28"

# zh-cn sheet: Priority (E) changed from "ht" to "mt" for rows 4-5
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback DateTime (K)
$wsZhCn.Range("H4").Value = "2016-08-22 18:16:23"
$wsZhCn.Range("H5").Value = "2016-08-22 18:16:23"
$wsZhCn.Range("K4").Value = "2016-08-22 18:16:39"
$wsZhCn.Range("K5").Value = "2016-08-22 18:16:39"

# de-de sheet: Priority (E) changed from "ht" to "mt" for rows 4-5
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (H)
$wsDeDe.Range("H4").Value = "2016-08-22 18:16:28"
$wsDeDe.Range("H5").Value = "2016-08-22 18:16:28"

# de-de sheet: Correspond Handback DateTime (K)
$wsDeDe.Range("K4").Value = "2016-08-22 18:16:46"
$wsDeDe.Range("K5").Value = "2016-08-22 18:16:46"
